$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before D, shifting existing D:K data to F:M
$ws.Range("D:E").EntireColumn.Insert()

# Row 7
$ws.Range("F7").Copy()
$ws.Range("D7:E7").PasteSpecial(-4122)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373

# Row 8
$ws.Range("F8").Copy()
$ws.Range("D8:E8").PasteSpecial(-4122)
$ws.Range("D8").Value = 4186800
$ws.Range("E8").Value = 3999400

# Row 9
$ws.Range("F9").Copy()
$ws.Range("D9:E9").PasteSpecial(-4122)
$ws.Range("D9").Value = 3069300
$ws.Range("E9").Value = 2412300

# Row 10
$ws.Range("F10").Copy()
$ws.Range("D10:E10").PasteSpecial(-4122)
$ws.Range("D10").Value = 1117600
$ws.Range("E10").Value = 1587100

# Row 11
$ws.Range("F11").Copy()
$ws.Range("D11:E11").PasteSpecial(-4122)

# Row 12
$ws.Range("F12").Copy()
$ws.Range("D12:E12").PasteSpecial(-4122)
$ws.Range("D12").Value = 276800
$ws.Range("E12").Value = 327000

# Row 13
$ws.Range("F13").Copy()
$ws.Range("D13:E13").PasteSpecial(-4122)
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0

# Row 14
$ws.Range("F14").Copy()
$ws.Range("D14:E14").PasteSpecial(-4122)
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0

# Row 15
$ws.Range("F15").Copy()
$ws.Range("D15:E15").PasteSpecial(-4122)
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0

# Row 16
$ws.Range("F16").Copy()
$ws.Range("D16:E16").PasteSpecial(-4122)

# Row 17
$ws.Range("F17").Copy()
$ws.Range("D17:E17").PasteSpecial(-4122)
$ws.Range("D17").Value = 3971100
$ws.Range("E17").Value = 3518700

# Row 18
$ws.Range("F18").Copy()
$ws.Range("D18:E18").PasteSpecial(-4122)
$ws.Range("D18").Value = 215700
$ws.Range("E18").Value = 480700

# Row 19
$ws.Range("F19").Copy()
$ws.Range("D19:E19").PasteSpecial(-4122)

# Row 20
$ws.Range("F20").Copy()
$ws.Range("D20:E20").PasteSpecial(-4122)
$ws.Range("D20").Value = 32500
$ws.Range("E20").Value = 7000

# Row 21
$ws.Range("F21").Copy()
$ws.Range("D21:E21").PasteSpecial(-4122)
$ws.Range("D21").Value = 2334000
$ws.Range("E21").Value = 2430600

# Row 22
$ws.Range("F22").Copy()
$ws.Range("D22:E22").PasteSpecial(-4122)
$ws.Range("D22").Value = 128800
$ws.Range("E22").Value = 108900

# Row 23
$ws.Range("F23").Copy()
$ws.Range("D23:E23").PasteSpecial(-4122)
$ws.Range("D23").Value = 119400
$ws.Range("E23").Value = 378800

# Row 24
$ws.Range("F24").Copy()
$ws.Range("D24:E24").PasteSpecial(-4122)
$ws.Range("D24").Value = -23200
$ws.Range("E24").Value = 13600

# Row 25
$ws.Range("F25").Copy()
$ws.Range("D25:E25").PasteSpecial(-4122)
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0

# Row 26
$ws.Range("F26").Copy()
$ws.Range("D26:E26").PasteSpecial(-4122)
$ws.Range("D26").Value = 142600
$ws.Range("E26").Value = 365200

# Row 27
$ws.Range("F27").Copy()
$ws.Range("D27:E27").PasteSpecial(-4122)
$ws.Range("D27").Value = 142600
$ws.Range("E27").Value = 365200

# Row 28
$ws.Range("F28").Copy()
$ws.Range("D28:E28").PasteSpecial(-4122)
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0

# Row 29
$ws.Range("F29").Copy()
$ws.Range("D29:E29").PasteSpecial(-4122)
$ws.Range("D29").Value = -8700
$ws.Range("E29").Value = 37600

# Row 30
$ws.Range("F30").Copy()
$ws.Range("D30:E30").PasteSpecial(-4122)
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0

# Row 31
$ws.Range("F31").Copy()
$ws.Range("D31:E31").PasteSpecial(-4122)
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0

# Row 32
$ws.Range("F32").Copy()
$ws.Range("D32:E32").PasteSpecial(-4122)
$ws.Range("D32").Value = -32500
$ws.Range("E32").Value = -7000

# Row 33
$ws.Range("F33").Copy()
$ws.Range("D33:E33").PasteSpecial(-4122)
$ws.Range("D33").Value = 133900
$ws.Range("E33").Value = 402800

# Row 34
$ws.Range("F34").Copy()
$ws.Range("D34:E34").PasteSpecial(-4122)
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0

# Row 35
$ws.Range("F35").Copy()
$ws.Range("D35:E35").PasteSpecial(-4122)
$ws.Range("D35").Value = 133900
$ws.Range("E35").Value = 402800

# Row 38
$ws.Range("F38").Copy()
$ws.Range("D38:E38").PasteSpecial(-4122)
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373

# Row 39
$ws.Range("F39").Copy()
$ws.Range("D39:E39").PasteSpecial(-4122)

# Row 40
$ws.Range("F40").Copy()
$ws.Range("D40:E40").PasteSpecial(-4122)

# Row 41
$ws.Range("F41").Copy()
$ws.Range("D41:E41").PasteSpecial(-4122)
$ws.Range("D41").Value = 3794500
$ws.Range("E41").Value = 3067500

# Row 42
$ws.Range("F42").Copy()
$ws.Range("D42:E42").PasteSpecial(-4122)
$ws.Range("D42").Value = "NA"
$ws.Range("E42").Value = "NA"

# Row 43
$ws.Range("F43").Copy()
$ws.Range("D43:E43").PasteSpecial(-4122)
$ws.Range("D43").Value = 5151200
$ws.Range("E43").Value = 4987900

# Row 44
$ws.Range("F44").Copy()
$ws.Range("D44:E44").PasteSpecial(-4122)
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0

# Row 45
$ws.Range("F45").Copy()
$ws.Range("D45:E45").PasteSpecial(-4122)
$ws.Range("D45").Value = 748500
$ws.Range("E45").Value = 674500

# Row 46
$ws.Range("F46").Copy()
$ws.Range("D46:E46").PasteSpecial(-4122)
$ws.Range("D46").Value = 9694100
$ws.Range("E46").Value = 8730000

# Row 47
$ws.Range("F47").Copy()
$ws.Range("D47:E47").PasteSpecial(-4122)
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0

# Row 48
$ws.Range("F48").Copy()
$ws.Range("D48:E48").PasteSpecial(-4122)
$ws.Range("D48").Value = 418300
$ws.Range("E48").Value = 371200

# Row 49
$ws.Range("F49").Copy()
$ws.Range("D49:E49").PasteSpecial(-4122)
$ws.Range("D49").Value = 14961000
$ws.Range("E49").Value = 13408400

# Row 50
$ws.Range("F50").Copy()
$ws.Range("D50:E50").PasteSpecial(-4122)
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0

# Row 51
$ws.Range("F51").Copy()
$ws.Range("D51:E51").PasteSpecial(-4122)
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0

# Row 52
$ws.Range("F52").Copy()
$ws.Range("D52:E52").PasteSpecial(-4122)
$ws.Range("D52").Value = 901000
$ws.Range("E52").Value = 856700

# Row 53
$ws.Range("F53").Copy()
$ws.Range("D53:E53").PasteSpecial(-4122)
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0

# Row 54
$ws.Range("F54").Copy()
$ws.Range("D54:E54").PasteSpecial(-4122)
$ws.Range("D54").Value = 25974400
$ws.Range("E54").Value = 23366200

# Row 55
$ws.Range("F55").Copy()
$ws.Range("D55:E55").PasteSpecial(-4122)

# Row 56
$ws.Range("F56").Copy()
$ws.Range("D56:E56").PasteSpecial(-4122)

# Row 57
$ws.Range("F57").Copy()
$ws.Range("D57:E57").PasteSpecial(-4122)
$ws.Range("D57").Value = 5249000
$ws.Range("E57").Value = 5054400

# Row 58
$ws.Range("F58").Copy()
$ws.Range("D58:E58").PasteSpecial(-4122)
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0

# Row 59
$ws.Range("F59").Copy()
$ws.Range("D59:E59").PasteSpecial(-4122)
$ws.Range("D59").Value = 1238300
$ws.Range("E59").Value = 1243800

# Row 60
$ws.Range("F60").Copy()
$ws.Range("D60:E60").PasteSpecial(-4122)
$ws.Range("D60").Value = 6487300
$ws.Range("E60").Value = 6298200

# Row 61
$ws.Range("F61").Copy()
$ws.Range("D61:E61").PasteSpecial(-4122)
$ws.Range("D61").Value = 10360100
$ws.Range("E61").Value = 8336600

# Row 62
$ws.Range("F62").Copy()
$ws.Range("D62:E62").PasteSpecial(-4122)
$ws.Range("D62").Value = 3888300
$ws.Range("E62").Value = 3721800

# Row 63
$ws.Range("F63").Copy()
$ws.Range("D63:E63").PasteSpecial(-4122)
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0

# Row 64
$ws.Range("F64").Copy()
$ws.Range("D64:E64").PasteSpecial(-4122)
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0

# Row 65
$ws.Range("F65").Copy()
$ws.Range("D65:E65").PasteSpecial(-4122)
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0

# Row 66
$ws.Range("F66").Copy()
$ws.Range("D66:E66").PasteSpecial(-4122)
$ws.Range("D66").Value = 20735600
$ws.Range("E66").Value = 18356600

# Row 67
$ws.Range("F67").Copy()
$ws.Range("D67:E67").PasteSpecial(-4122)

# Row 68
$ws.Range("F68").Copy()
$ws.Range("D68:E68").PasteSpecial(-4122)
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0

# Row 69
$ws.Range("F69").Copy()
$ws.Range("D69:E69").PasteSpecial(-4122)
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0

# Row 70
$ws.Range("F70").Copy()
$ws.Range("D70:E70").PasteSpecial(-4122)
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0

# Row 71
$ws.Range("F71").Copy()
$ws.Range("D71:E71").PasteSpecial(-4122)
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0

# Row 72
$ws.Range("F72").Copy()
$ws.Range("D72:E72").PasteSpecial(-4122)
$ws.Range("D72").Value = 2942400
$ws.Range("E72").Value = 2808400

# Row 73
$ws.Range("F73").Copy()
$ws.Range("D73:E73").PasteSpecial(-4122)
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0

# Row 74
$ws.Range("F74").Copy()
$ws.Range("D74:E74").PasteSpecial(-4122)
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0

# Row 75
$ws.Range("F75").Copy()
$ws.Range("D75:E75").PasteSpecial(-4122)
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0

# Row 76
$ws.Range("F76").Copy()
$ws.Range("D76:E76").PasteSpecial(-4122)
$ws.Range("D76").Value = 5238800
$ws.Range("E76").Value = 5009700

# Row 77
$ws.Range("F77").Copy()
$ws.Range("D77:E77").PasteSpecial(-4122)
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0

# Row 80
$ws.Range("F80").Copy()
$ws.Range("D80:E80").PasteSpecial(-4122)
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373

# Row 81
$ws.Range("F81").Copy()
$ws.Range("D81:E81").PasteSpecial(-4122)
$ws.Range("D81").Value = 133900
$ws.Range("E81").Value = 402800

# Row 82
$ws.Range("F82").Copy()
$ws.Range("D82:E82").PasteSpecial(-4122)

# Row 83
$ws.Range("F83").Copy()
$ws.Range("D83:E83").PasteSpecial(-4122)
$ws.Range("D83").Value = 2076900
$ws.Range("E83").Value = 1932900

# Row 84
$ws.Range("F84").Copy()
$ws.Range("D84:E84").PasteSpecial(-4122)
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0

# Row 85
$ws.Range("F85").Copy()
$ws.Range("D85:E85").PasteSpecial(-4122)
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0

# Row 86
$ws.Range("F86").Copy()
$ws.Range("D86:E86").PasteSpecial(-4122)
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0

# Row 87
$ws.Range("F87").Copy()
$ws.Range("D87:E87").PasteSpecial(-4122)
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0

# Row 88
$ws.Range("F88").Copy()
$ws.Range("D88:E88").PasteSpecial(-4122)
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0

# Row 89
$ws.Range("F89").Copy()
$ws.Range("D89:E89").PasteSpecial(-4122)
$ws.Range("D89").Value = -1235100
$ws.Range("E89").Value = -690400

# Row 90
$ws.Range("F90").Copy()
$ws.Range("D90:E90").PasteSpecial(-4122)

# Row 91
$ws.Range("F91").Copy()
$ws.Range("D91:E91").PasteSpecial(-4122)
$ws.Range("D91").Value = -70100
$ws.Range("E91").Value = -39300

# Row 92
$ws.Range("F92").Copy()
$ws.Range("D92:E92").PasteSpecial(-4122)
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0

# Row 93
$ws.Range("F93").Copy()
$ws.Range("D93:E93").PasteSpecial(-4122)
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0

# Row 94
$ws.Range("F94").Copy()
$ws.Range("D94:E94").PasteSpecial(-4122)
$ws.Range("D94").Value = -80400
$ws.Range("E94").Value = -168700

# Row 95
$ws.Range("F95").Copy()
$ws.Range("D95:E95").PasteSpecial(-4122)

# Row 96
$ws.Range("F96").Copy()
$ws.Range("D96:E96").PasteSpecial(-4122)
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0

# Row 97
$ws.Range("F97").Copy()
$ws.Range("D97:E97").PasteSpecial(-4122)
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0

# Row 98
$ws.Range("F98").Copy()
$ws.Range("D98:E98").PasteSpecial(-4122)
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0

# Row 99
$ws.Range("F99").Copy()
$ws.Range("D99:E99").PasteSpecial(-4122)
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0

# Row 100
$ws.Range("F100").Copy()
$ws.Range("D100:E100").PasteSpecial(-4122)
$ws.Range("D100").Value = 2053900
$ws.Range("E100").Value = 29200

# Row 101
$ws.Range("F101").Copy()
$ws.Range("D101:E101").PasteSpecial(-4122)
$ws.Range("D101").Value = -5000
$ws.Range("E101").Value = -5600

# Row 102
$ws.Range("F102").Copy()
$ws.Range("D102:E102").PasteSpecial(-4122)
$ws.Range("D102").Value = 733500
$ws.Range("E102").Value = -835400

# Row 91 special case: F:J values are restated (not a simple shift of old D:H)
$ws.Range("F91").Value = -27300
$ws.Range("G91").Value = -37200
$ws.Range("H91").Value = -21600
$ws.Range("I91").Value = -34000
$ws.Range("J91").Value = -65200
